$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Allosaurus Rider', ['{5}{G}{G}', 'Creature — Elf Warrior', 'You may exile two green cards from your hand rather than pay this spell’s mana cost.', 'Allosaurus Rider’s power and toughness are each equal to 1 plus the number of lands you control.', '1+*/1+*'])"
$ws.Range("A3").Value = "('Marit Lage', ['Token Legendary Creature — Avatar', 'Flying, indestructible', '20/20'])"

$ws.Range("A4:A11").EntireRow.Delete()
